# "new code testng - dataprovider": populate the Result column (Pass/Fail)
# used by the TestNG DataProvider-driven login test on the Credentials sheet,
# then leave the selection where the author's editor session left it.

$wb = $excel.ActiveWorkbook

# --- Credentials sheet: add Result column data (Pass/Fail) ---
$wsCred = $wb.Worksheets.Item("Credentials")
$wsCred.Range("C2").Value = "Pass"   # standard_user -> Pass
$wsCred.Range("C3").Value = "Fail"   # locked_out_user -> Fail
$wsCred.Range("C4").Value = "Pass"   # problem_user -> Pass

# --- Sheet1: move selection back to A1 ---
$wsSheet1 = $wb.Worksheets.Item("Sheet1")
$wsSheet1.Range("A1").Select()

# --- Reselect Credentials as the active/tab-selected sheet with A2 active cell ---
$wsCred.Activate()
$wsCred.Range("A2").Select()
